$d = $word.ActiveDocument

$replacements = @(
    @{old = "91÷9=10, 1"; new = "93÷2=46, 1"},
    @{old = "95÷3=31, 2"; new = "25÷4=6, 1"},
    @{old = "70÷8=8, 6"; new = "97÷7=13, 6"},
    @{old = "56÷7=8, 0"; new = "53÷8=6, 5"},
    @{old = "41÷4=10, 1"; new = "22÷6=3, 4"},
    @{old = "50÷2=25, 0"; new = "56÷2=28, 0"},
    @{old = "99÷7=14, 1"; new = "86÷9=9, 5"},
    @{old = "53÷4=13, 1"; new = "26÷3=8, 2"},
    @{old = "90÷8=11, 2"; new = "61÷4=15, 1"},
    @{old = "50÷3=16, 2"; new = "92÷6=15, 2"},
    @{old = "74÷8=9, 2"; new = "23÷7=3, 2"},
    @{old = "83÷3=27, 2"; new = "99÷5=19, 4"},
    @{old = "88÷8=11, 0"; new = "51÷4=12, 3"},
    @{old = "11÷5=2, 1"; new = "44÷9=4, 8"},
    @{old = "47÷8=5, 7"; new = "82÷2=41, 0"},
    @{old = "68÷8=8, 4"; new = "37÷7=5, 2"},
    @{old = "79÷7=11, 2"; new = "68÷9=7, 5"},
    @{old = "11÷4=2, 3"; new = "31÷5=6, 1"},
    @{old = "72÷4=18, 0"; new = "66÷7=9, 3"},
    @{old = "72÷5=14, 2"; new = "69÷2=34, 1"},
    @{old = "40÷2=20, 0"; new = "36÷4=9, 0"},
    @{old = "71÷9=7, 8"; new = "17÷4=4, 1"},
    @{old = "24÷5=4, 4"; new = "44÷2=22, 0"},
    @{old = "20÷6=3, 2"; new = "97÷3=32, 1"},
    @{old = "89÷7=12, 5"; new = "17÷5=3, 2"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
